$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Clear the placeholder zero values in row 1, columns C:F.
# ClearContents() removes the stored <v>0</v> values but keeps the
# existing cell formatting/style (s="2") intact, matching the diff.
$ws.Range("C1:F1").ClearContents()

# Move/record the active selection on the sheet to H6.
$ws.Range("H6").Select()
